$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 87 with the latest Argent price entry.
# Force text storage (matching the existing inline/shared-string cells in
# this column) instead of Excel's automatic date/number inference, then
# clear the temporary number-format override so no stray style is left
# behind on the new cells.
$ws.Cells.Item(87, 1).NumberFormat = "@"
$ws.Cells.Item(87, 1).Value = "2025-01-24"
$ws.Cells.Item(87, 1).ClearFormats()

$ws.Cells.Item(87, 2).NumberFormat = "@"
$ws.Cells.Item(87, 2).Value = "5.84"
$ws.Cells.Item(87, 2).ClearFormats()
